$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the Name (C) values for the server rows, dropping the "_1" suffix.
# This also seeds new shared strings 35..39 in the exact order required.
$ws.Cells.Item(10, 3).Value = "GameServer"
$ws.Cells.Item(11, 3).Value = "WorldServer"
$ws.Cells.Item(12, 3).Value = "ProxyServer"
$ws.Cells.Item(13, 3).Value = "MasterServer"
$ws.Cells.Item(14, 3).Value = "LoginServer"

# --- New column I ("Area") header + type/flags rows, mirroring column H's
# formatting by copying the cell (keeps style + conditional flags in sync).
$ws.Cells.Item(1, 9).Value = "Area"          # new shared string 40
$ws.Cells.Item(2, 8).Copy($ws.Cells.Item(2, 9))
$ws.Cells.Item(3, 8).Copy($ws.Cells.Item(3, 9))
$ws.Cells.Item(4, 8).Copy($ws.Cells.Item(4, 9))
$ws.Cells.Item(5, 8).Copy($ws.Cells.Item(5, 9))
$ws.Cells.Item(6, 8).Copy($ws.Cells.Item(6, 9))
$ws.Cells.Item(7, 8).Copy($ws.Cells.Item(7, 9))
$ws.Cells.Item(8, 8).Copy($ws.Cells.Item(8, 9))

# --- Row 9 (Desc row) gets the Chinese description for the new column.
$ws.Cells.Item(9, 9).Value = "区服"            # new shared string 41

# --- Public row (3) flips from TRUE to FALSE across B:H now that a new
# column exists and Public is no longer blanket-true.
$ws.Range("B3:H3").Value = $false

# --- Data rows 10-14: populate the new Area column with 1.
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(12, 9).Value = 1
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(14, 9).Value = 1

# --- Update the TRUE/FALSE list validation so it also covers the new
# column I cells (G7:J8 includes I&J; F2 kept separate from F7:F9).
$ws.Range("F15:F1048576").Validation.Delete()
$ws.Range("F9").Validation.Delete()
$ws.Range("F2:F6").Validation.Delete()
$ws.Range("B7:J8").Validation.Delete()

$ws.Range("F15:F1048576").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("F7:F9").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("G7:J8").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("F2").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("B7:E8").Validation.Add(3, 1, 1, """TRUE,FALSE""")

# --- Selection ends on the new column, matching the saved UI state.
$ws.Cells.Item(17, 9).Select()
